$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C26").Value = 3
$ws.Range("E26").Value = 160

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.HasTitle = $false
